# H4 nagelezen voor een eerste keer
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Range" query-type column header
$ws.Range("G2").Value = "Range"

# Row 3 (Cassandra / Column)
$ws.Range("G3").Value = "Half"

# Row 5 (Document / Apache CouchDB)
$ws.Range("F5").Value = "Nee"
$ws.Range("G5").Value = "Ja"

# Row 8 (Key-Value / LightCloud) - fill in previously empty cells
$ws.Range("E8").Value = "Ja"
$ws.Range("F8").Value = "Nee"
$ws.Range("G8").Value = "Ja"

# Row 10 (MemcacheDB)
$ws.Range("D10").Value = "Ja"
$ws.Range("E10").Value = "Nee"
$ws.Range("F10").Value = "Nee"
$ws.Range("G10").Value = "Ja"

# Row 11 (Redis)
$ws.Range("C11").Value = "Snapshots"
$ws.Range("D11").Value = "Master-Slave"
$ws.Range("D11").VerticalAlignment = -4108
$ws.Range("E11").Value = "Half"
$ws.Range("F11").Value = "Ja"
$ws.Range("G11").Value = "Half"

# Row 12 (Riak)
$ws.Range("E12").Value = "Ja"
$ws.Range("F12").Value = "Nee"
$ws.Range("G12").Value = "Half"

# Row 13 (Voldemort)
$ws.Range("E13").Value = "Ja"
$ws.Range("F13").Value = "Nee"
$ws.Range("G13").Value = "Nee"

# Row 14 (Relationeel / MySQL)
$ws.Range("D14").Value = "Master-Slave"
$ws.Range("E14").Value = "Nee"
$ws.Range("F14").Value = "Ja"
$ws.Range("G14").Value = "Ja"

# Row 15 (Pgpool-II)
$ws.Range("E15").Value = "Mogelijk"
$ws.Range("F15").Value = "Ja"
$ws.Range("G15").Value = "Ja"

# New column width for column F (added when the "Range" column content was entered)
$ws.Columns.Item(6).ColumnWidth = 9.6

# Update the last active-cell selection recorded in the sheet view
$ws.Range("O12").Select()
